# Update cryptocurrency Price (column D) and Volume(1h) (column E) figures
# with the latest scraped values (GitHub Actions refresh, 2023-01-19 08:48 UTC).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").Value = "'291.35"
$ws.Range("E2").Value = "'-3.28%"

# Row 3: OKB
$ws.Range("D3").Value = "'30.66"
$ws.Range("E3").Value = "'-6.55%"

# Row 4: HuobiToken
$ws.Range("D4").Value = "'4.952"
$ws.Range("E4").Value = "'0.09%"

# Row 5: Cronos
$ws.Range("D5").Value = "'0.07238"
$ws.Range("E5").Value = "'-6.39%"

# Row 6: FTXToken
$ws.Range("D6").Value = "'1.793"
$ws.Range("E6").Value = "'-7.82%"

# Row 7: KuCoinToken
$ws.Range("D7").Value = "'7.688"
$ws.Range("E7").Value = "'-1.89%"

# Row 8: GateToken
$ws.Range("D8").Value = "'3.756"
$ws.Range("E8").Value = "'-1.25%"

# Row 9: MXToken
$ws.Range("D9").Value = "'0.8990"
$ws.Range("E9").Value = "'-2.31%"

# Row 10: WazirX
$ws.Range("D10").Value = "'0.1661"
$ws.Range("E10").Value = "'-5.98%"

# Row 11: LiechtensteinCryptoassetsExchange
$ws.Range("D11").Value = "'0.07707"
$ws.Range("E11").Value = "'-1.67%"

# Row 12: MandalaExchangeToken
$ws.Range("D12").Value = "'0.08029"
$ws.Range("E12").Value = "'-7.17%"

# Row 13: BitrueCoin
$ws.Range("D13").Value = "'0.03034"
$ws.Range("E13").Value = "'-4.27%"

# Row 14: BitMartToken
$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'-0.18%"

# Row 15: BitForexToken
$ws.Range("D15").Value = "'0.001496"
$ws.Range("E15").Value = "'-1.21%"

# Row 16: TigerCash
$ws.Range("D16").Value = "'0.005691"
$ws.Range("E16").Value = "'-1.96%"

# Row 18: LEO
$ws.Range("D18").Value = "'3.472"
$ws.Range("E18").Value = "'0.29%"

# Row 19: BTSEToken
$ws.Range("E19").Value = "'-3.24%"

# Row 20: BitpandaEcosystemToken
$ws.Range("E20").Value = "'-0.86%"

# Row 21: ProBitToken
$ws.Range("D21").Value = "'0.1288"
$ws.Range("E21").Value = "'-2.93%"

# Row 22: MCDex
$ws.Range("D22").Value = "'4.046"
$ws.Range("E22").Value = "'-5.63%"

# Row 23: ZBToken
$ws.Range("D23").Value = "'0.2250"
$ws.Range("E23").Value = "'12.88%"

# Row 24: CoinExToken
$ws.Range("D24").Value = "'0.04500"
$ws.Range("E24").Value = "'-1.06%"

# Row 25: BitKan
$ws.Range("D25").Value = "'0.001216"
$ws.Range("E25").Value = "'-0.83%"

# Row 26: HotbitToken
$ws.Range("D26").Value = "'0.004017"
$ws.Range("E26").Value = "'-8.95%"

# Row 27: NitroEx
$ws.Range("D27").Value = "'0.0001251"
$ws.Range("E27").Value = "'-0.04%"

# Row 39: One
$ws.Range("D39").Value = "'0.01608"
$ws.Range("E39").Value = "'-5.76%"

# Row 40: IDEX
$ws.Range("D40").Value = "'0.04420"
$ws.Range("E40").Value = "'-5.69%"

# Row 41: KickToken
$ws.Range("D41").Value = "'0.007284"
$ws.Range("E41").Value = "'-5.18%"

# Row 42: BKEXToken
$ws.Range("D42").Value = "'0.1309"
$ws.Range("E42").Value = "'-3.05%"

# Row 43: Dexo
$ws.Range("D43").Value = "'0.007755"

# Row 44: CEJI
$ws.Range("D44").Value = "'0.002015"
$ws.Range("E44").Value = "'-13.28%"

# Row 45: LocalTraders
$ws.Range("D45").Value = "'0.009510"
$ws.Range("E45").Value = "'-16.46%"

# Row 46: CoinLion
$ws.Range("E46").Value = "'-5.31%"

# Row 47: Kangarootoken
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.17%"

# Row 48: BOLO
$ws.Range("D48").Value = "'2.247"
$ws.Range("E48").Value = "'173.92%"

# Row 49: CoinbaseStockToken
$ws.Range("D49").Value = "'0.002999"
$ws.Range("E49").Value = "'-3.41%"

# Row 50: CryptobidCoin
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.17%"

# Row 51: SpecialPowerGold
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.17%"
